$wb = $excel.ActiveWorkbook

# Update the "想去人数" (attendance count) figures for two rows on both the
# "展览" (Exhibition) and "全部类型" (All types) sheets, which carry
# duplicated data.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 462
    $ws.Range("F4").Value = 22
}
